# Updated cryptos list on Wed Feb 28 13:30:19 UTC 2024 with GitHub Actions
#
# Refreshes the Price (column D) and Volume(1h) (column E) columns for
# rows 2-51 of the crypto tracker sheet.
#
# Note: several Price values look like plain numbers (e.g. "0.998",
# "3.35", "1.00"). Excel's Range.Value setter auto-detects numeric-looking
# text and stores it as a Number, which would silently reformat values
# such as "1.00" -> "1" and drop the original text representation. The
# source data keeps these as literal text, so a leading apostrophe is used
# to force text entry (standard Excel convention for "treat this as text"),
# exactly like typing '0.998 into a cell. The apostrophe itself is not
# stored; Excel just marks the cell as text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.298.94"
$ws.Range("E2").Value = "  +5.88%  "

$ws.Range("D3").Value = "3.347.09"
$ws.Range("E3").Value = "  +2.54%  "

$ws.Range("D4").Value = "'0.998"
$ws.Range("E4").Value = "  -0.33%  "

$ws.Range("D5").Value = "'411.76"
$ws.Range("E5").Value = "  +3.38%  "

$ws.Range("D6").Value = "'111.86"
$ws.Range("E6").Value = "  +1.94%  "

$ws.Range("D7").Value = "'0.584"
$ws.Range("E7").Value = "  +4.92%  "

$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = "  -0.16%  "

$ws.Range("D9").Value = "'0.633"
$ws.Range("E9").Value = "  +1.96%  "

$ws.Range("D10").Value = "'40.02"
$ws.Range("E10").Value = "  +1.98%  "

$ws.Range("D11").Value = "'0.0989"
$ws.Range("E11").Value = "  +3.84%  "

$ws.Range("D12").Value = "'0.144"
$ws.Range("E12").Value = "  +1.34%  "

$ws.Range("D13").Value = "3.860.28"
$ws.Range("E13").Value = "  +1.88%  "

$ws.Range("D14").Value = "'8.55"
$ws.Range("E14").Value = "  +6.08%  "

$ws.Range("D15").Value = "'19.34"
$ws.Range("E15").Value = "  +1.55%  "

$ws.Range("D16").Value = "3.344.08"
$ws.Range("E16").Value = "  +2.47%  "

$ws.Range("D17").Value = "'1.05"
$ws.Range("E17").Value = "  +0.49%  "

$ws.Range("D18").Value = "59.888.72"
$ws.Range("E18").Value = "  +5.27%  "

$ws.Range("D19").Value = "'10.69"
$ws.Range("E19").Value = "  -0.97%  "

$ws.Range("D20").Value = "'3.35"
$ws.Range("E20").Value = "  +1.59%  "

$ws.Range("D21").Value = "'0.0000111"
$ws.Range("E21").Value = "  +4.65%  "

$ws.Range("D22").Value = "'13.09"
$ws.Range("E22").Value = "  +1.80%  "

$ws.Range("D23").Value = "'304.90"
$ws.Range("E23").Value = "  +0.30%  "

$ws.Range("D24").Value = "'75.63"
$ws.Range("E24").Value = "  +0.90%  "

$ws.Range("D25").Value = "'3.29"
$ws.Range("E25").Value = "  +3.46%  "

$ws.Range("D26").Value = "'0.185"
$ws.Range("E26").Value = "  +9.43%  "

$ws.Range("D27").Value = "'28.64"
$ws.Range("E27").Value = "  +1.50%  "

$ws.Range("D28").Value = "'4.47"
$ws.Range("E28").Value = "  +2.06%  "

$ws.Range("D29").Value = "'7.87"
$ws.Range("E29").Value = "  -1.08%  "

$ws.Range("D30").Value = "'7.47"
$ws.Range("E30").Value = "  +2.97%  "

$ws.Range("D31").Value = "'2.66"
$ws.Range("E31").Value = "  +24.90%  "

$ws.Range("E32").Value = "  +4.14%  "

$ws.Range("E33").Value = "  +0.05%  "

$ws.Range("D34").Value = "'11.56"
$ws.Range("E34").Value = "  +4.99%  "

$ws.Range("D35").Value = "'39.70"
$ws.Range("E35").Value = "  +6.47%  "

$ws.Range("D36").Value = "'0.0509"
$ws.Range("E36").Value = "  +5.30%  "

$ws.Range("D37").Value = "'51.92"
$ws.Range("E37").Value = "  +0.80%  "

$ws.Range("D38").Value = "'3.22"
$ws.Range("E38").Value = "  +1.43%  "

$ws.Range("D39").Value = "'0.996"
$ws.Range("E39").Value = "  -0.48%  "

$ws.Range("D40").Value = "'3.40"
$ws.Range("E40").Value = "  -4.33%  "

$ws.Range("D41").Value = "'138.70"
$ws.Range("E41").Value = "  +4.06%  "

$ws.Range("E42").Value = "  +2.75%  "

$ws.Range("E43").Value = "  -0.38%  "

$ws.Range("D44").Value = "'0.286"
$ws.Range("E44").Value = "  +1.91%  "

$ws.Range("D45").Value = "'3.94"
$ws.Range("E45").Value = "  -0.86%  "

$ws.Range("D46").Value = "'16.88"
$ws.Range("E46").Value = "  -2.92%  "

$ws.Range("E47").Value = "  +8.88%  "

$ws.Range("D48").Value = "'22.39"
$ws.Range("E48").Value = "  +2.12%  "

$ws.Range("D49").Value = "2.202.70"
$ws.Range("E49").Value = "  +2.32%  "

$ws.Range("D50").Value = "'2.10"
$ws.Range("E50").Value = "  +3.64%  "

$ws.Range("E51").Value = "  +0.16%  "
